$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 135,2
$arr[0,0] = 0.8466431206713463
$arr[0,1] = 0.1533568793286537
$arr[1,0] = 0.7820840370477621
$arr[1,1] = 0.2179159629522379
$arr[2,0] = 0.8749414277318412
$arr[2,1] = 0.1250585722681588
$arr[3,0] = 0.7820840370477621
$arr[3,1] = 0.2179159629522379
$arr[4,0] = 0.7602848156868673
$arr[4,1] = 0.2397151843131327
$arr[5,0] = 0.6420799356935972
$arr[5,1] = 0.3579200643064028
$arr[6,0] = 0.7820840370477621
$arr[6,1] = 0.2179159629522379
$arr[7,0] = 0.8466431206713463
$arr[7,1] = 0.1533568793286537
$arr[8,0] = 0.5659720316352707
$arr[8,1] = 0.4340279683647293
$arr[9,0] = 0.8994915014777216
$arr[9,1] = 0.1005084985222784
$arr[10,0] = 0.8466431206713463
$arr[10,1] = 0.1533568793286537
$arr[11,0] = 0.7820840370477621
$arr[11,1] = 0.2179159629522379
$arr[12,0] = 0.8994915014777216
$arr[12,1] = 0.1005084985222784
$arr[13,0] = 0.8994915014777216
$arr[13,1] = 0.1005084985222784
$arr[14,0] = 0.8994915014777216
$arr[14,1] = 0.1005084985222784
$arr[15,0] = 0.835482401800504
$arr[15,1] = 0.164517598199496
$arr[16,0] = 0.8994915014777216
$arr[16,1] = 0.1005084985222784
$arr[17,0] = 0.5498452523359121
$arr[17,1] = 0.4501547476640879
$arr[18,0] = 0.8466431206713463
$arr[18,1] = 0.1533568793286537
$arr[19,0] = 0.7820840370477621
$arr[19,1] = 0.2179159629522379
$arr[20,0] = 0.8749414277318412
$arr[20,1] = 0.1250585722681588
$arr[21,0] = 0.7329069162497991
$arr[21,1] = 0.2670930837502009
$arr[22,0] = 0.8466431206713463
$arr[22,1] = 0.1533568793286537
$arr[23,0] = 0.7820840370477621
$arr[23,1] = 0.2179159629522379
$arr[24,0] = 0.8994915014777216
$arr[24,1] = 0.1005084985222784
$arr[25,0] = 0.6401618593643882
$arr[25,1] = 0.3598381406356118
$arr[26,0] = 0.8466431206713463
$arr[26,1] = 0.1533568793286537
$arr[27,0] = 0.7820840370477621
$arr[27,1] = 0.2179159629522379
$arr[28,0] = 0.8466431206713463
$arr[28,1] = 0.1533568793286537
$arr[29,0] = 0.7820840370477621
$arr[29,1] = 0.2179159629522379
$arr[30,0] = 0.8994915014777216
$arr[30,1] = 0.1005084985222784
$arr[31,0] = 0.8466431206713463
$arr[31,1] = 0.1533568793286537
$arr[32,0] = 0.8749414277318412
$arr[32,1] = 0.1250585722681588
$arr[33,0] = 0.8994915014777216
$arr[33,1] = 0.1005084985222784
$arr[34,0] = 0.863300632083034
$arr[34,1] = 0.136699367916966
$arr[35,0] = 0.6335946406798525
$arr[35,1] = 0.3664053593201475
$arr[36,0] = 0.8466431206713463
$arr[36,1] = 0.1533568793286537
$arr[37,0] = 0.6661742215619185
$arr[37,1] = 0.3338257784380815
$arr[38,0] = 0.6537986360375384
$arr[38,1] = 0.3462013639624616
$arr[39,0] = 0.8994915014777216
$arr[39,1] = 0.1005084985222784
$arr[40,0] = 0.5998153575678875
$arr[40,1] = 0.4001846424321125
$arr[41,0] = 0.5498452523359121
$arr[41,1] = 0.4501547476640879
$arr[42,0] = 0.7607555257210576
$arr[42,1] = 0.2392444742789424
$arr[43,0] = 0.4829102216903616
$arr[43,1] = 0.5170897783096384
$arr[44,0] = 0.801630212515808
$arr[44,1] = 0.198369787484192
$arr[45,0] = 0.5071022729512122
$arr[45,1] = 0.4928977270487878
$arr[46,0] = 0.8994915014777216
$arr[46,1] = 0.1005084985222784
$arr[47,0] = 0.8749414277318412
$arr[47,1] = 0.1250585722681588
$arr[48,0] = 0.7820840370477621
$arr[48,1] = 0.2179159629522379
$arr[49,0] = 0.8466431206713463
$arr[49,1] = 0.1533568793286537
$arr[50,0] = 0.8466431206713463
$arr[50,1] = 0.1533568793286537
$arr[51,0] = 0.7820840370477621
$arr[51,1] = 0.2179159629522379
$arr[52,0] = 0.7820840370477621
$arr[52,1] = 0.2179159629522379
$arr[53,0] = 0.8749414277318412
$arr[53,1] = 0.1250585722681588
$arr[54,0] = 0.4309734809161656
$arr[54,1] = 0.5690265190838344
$arr[55,0] = 0.8994915014777216
$arr[55,1] = 0.1005084985222784
$arr[56,0] = 0.8700720494148051
$arr[56,1] = 0.1299279505851949
$arr[57,0] = 0.7820840370477621
$arr[57,1] = 0.2179159629522379
$arr[58,0] = 0.6612729355809915
$arr[58,1] = 0.3387270644190085
$arr[59,0] = 0.8994915014777216
$arr[59,1] = 0.1005084985222784
$arr[60,0] = 0.8994915014777216
$arr[60,1] = 0.1005084985222784
$arr[61,0] = 0.8994915014777216
$arr[61,1] = 0.1005084985222784
$arr[62,0] = 0.8954638149707839
$arr[62,1] = 0.1045361850292161
$arr[63,0] = 0.8466431206713463
$arr[63,1] = 0.1533568793286537
$arr[64,0] = 0.8466431206713463
$arr[64,1] = 0.1533568793286537
$arr[65,0] = 0.8700720494148051
$arr[65,1] = 0.1299279505851949
$arr[66,0] = 0.5998153575678875
$arr[66,1] = 0.4001846424321125
$arr[67,0] = 0.5655337270146484
$arr[67,1] = 0.4344662729853516
$arr[68,0] = 0.8466431206713463
$arr[68,1] = 0.1533568793286537
$arr[69,0] = 0.8994915014777216
$arr[69,1] = 0.1005084985222784
$arr[70,0] = 0.8466431206713463
$arr[70,1] = 0.1533568793286537
$arr[71,0] = 0.6467096640978485
$arr[71,1] = 0.3532903359021515
$arr[72,0] = 0.8994915014777216
$arr[72,1] = 0.1005084985222784
$arr[73,0] = 0.5464150277543774
$arr[73,1] = 0.4535849722456226
$arr[74,0] = 0.8749414277318412
$arr[74,1] = 0.1250585722681588
$arr[75,0] = 0.8466431206713463
$arr[75,1] = 0.1533568793286537
$arr[76,0] = 0.6181347653176266
$arr[76,1] = 0.3818652346823734
$arr[77,0] = 0.7820840370477621
$arr[77,1] = 0.2179159629522379
$arr[78,0] = 0.873418168661053
$arr[78,1] = 0.126581831338947
$arr[79,0] = 0.8994915014777216
$arr[79,1] = 0.1005084985222784
$arr[80,0] = 0.8466431206713463
$arr[80,1] = 0.1533568793286537
$arr[81,0] = 0.8994915014777216
$arr[81,1] = 0.1005084985222784
$arr[82,0] = 0.8466431206713463
$arr[82,1] = 0.1533568793286537
$arr[83,0] = 0.7820840370477621
$arr[83,1] = 0.2179159629522379
$arr[84,0] = 0.8749414277318412
$arr[84,1] = 0.1250585722681588
$arr[85,0] = 0.8466431206713463
$arr[85,1] = 0.1533568793286537
$arr[86,0] = 0.7820840370477621
$arr[86,1] = 0.2179159629522379
$arr[87,0] = 0.8856930230787471
$arr[87,1] = 0.1143069769212529
$arr[88,0] = 0.8101268523212063
$arr[88,1] = 0.1898731476787937
$arr[89,0] = 0.7441171837957787
$arr[89,1] = 0.2558828162042213
$arr[90,0] = 0.8171955475183693
$arr[90,1] = 0.1828044524816307
$arr[91,0] = 0.8994915014777216
$arr[91,1] = 0.1005084985222784
$arr[92,0] = 0.8749414277318412
$arr[92,1] = 0.1250585722681588
$arr[93,0] = 0.8994915014777216
$arr[93,1] = 0.1005084985222784
$arr[94,0] = 0.8466431206713463
$arr[94,1] = 0.1533568793286537
$arr[95,0] = 0.7797145041502886
$arr[95,1] = 0.2202854958497114
$arr[96,0] = 0.6406594588827165
$arr[96,1] = 0.3593405411172835
$arr[97,0] = 0.801630212515808
$arr[97,1] = 0.198369787484192
$arr[98,0] = 0.8749414277318412
$arr[98,1] = 0.1250585722681588
$arr[99,0] = 0.8749414277318412
$arr[99,1] = 0.1250585722681588
$arr[100,0] = 0.7031997836794691
$arr[100,1] = 0.2968002163205309
$arr[101,0] = 0.7820840370477621
$arr[101,1] = 0.2179159629522379
$arr[102,0] = 0.8466431206713463
$arr[102,1] = 0.1533568793286537
$arr[103,0] = 0.6420799356935972
$arr[103,1] = 0.3579200643064028
$arr[104,0] = 0.6420799356935972
$arr[104,1] = 0.3579200643064028
$arr[105,0] = 0.8466431206713463
$arr[105,1] = 0.1533568793286537
$arr[106,0] = 0.8466431206713463
$arr[106,1] = 0.1533568793286537
$arr[107,0] = 0.8466431206713463
$arr[107,1] = 0.1533568793286537
$arr[108,0] = 0.8749414277318412
$arr[108,1] = 0.1250585722681588
$arr[109,0] = 0.6661742215619185
$arr[109,1] = 0.3338257784380815
$arr[110,0] = 0.8380566359458733
$arr[110,1] = 0.1619433640541267
$arr[111,0] = 0.7820840370477621
$arr[111,1] = 0.2179159629522379
$arr[112,0] = 0.824948703925055
$arr[112,1] = 0.175051296074945
$arr[113,0] = 0.8466431206713463
$arr[113,1] = 0.1533568793286537
$arr[114,0] = 0.8466431206713463
$arr[114,1] = 0.1533568793286537
$arr[115,0] = 0.8466431206713463
$arr[115,1] = 0.1533568793286537
$arr[116,0] = 0.6590750771611322
$arr[116,1] = 0.3409249228388678
$arr[117,0] = 0.8749414277318412
$arr[117,1] = 0.1250585722681588
$arr[118,0] = 0.8466431206713463
$arr[118,1] = 0.1533568793286537
$arr[119,0] = 0.8466431206713463
$arr[119,1] = 0.1533568793286537
$arr[120,0] = 0.8749414277318412
$arr[120,1] = 0.1250585722681588
$arr[121,0] = 0.5960669561043807
$arr[121,1] = 0.4039330438956193
$arr[122,0] = 0.8749414277318412
$arr[122,1] = 0.1250585722681588
$arr[123,0] = 0.8466431206713463
$arr[123,1] = 0.1533568793286537
$arr[124,0] = 0.8749414277318412
$arr[124,1] = 0.1250585722681588
$arr[125,0] = 0.8994915014777216
$arr[125,1] = 0.1005084985222784
$arr[126,0] = 0.8994915014777216
$arr[126,1] = 0.1005084985222784
$arr[127,0] = 0.8466431206713463
$arr[127,1] = 0.1533568793286537
$arr[128,0] = 0.8466431206713463
$arr[128,1] = 0.1533568793286537
$arr[129,0] = 0.8380566359458733
$arr[129,1] = 0.1619433640541267
$arr[130,0] = 0.863300632083034
$arr[130,1] = 0.136699367916966
$arr[131,0] = 0.7042392952545217
$arr[131,1] = 0.2957607047454783
$arr[132,0] = 0.8003118571500073
$arr[132,1] = 0.1996881428499927
$arr[133,0] = 0.7602848156868673
$arr[133,1] = 0.2397151843131327
$arr[134,0] = 0.7820840370477621
$arr[134,1] = 0.2179159629522379
$ws.Range("A2:B136").Value = $arr
